# Insert a new record row at row 73 (Femacal de La Calera / Ciboulette),
# which pushes the existing rows 73..357 down to 74..358 and extends the
# table by one row (dimension A1:R357 -> A1:R358).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(73).Insert()

$ws.Range("A73").Value = 3
$ws.Range("B73").Value = "Femacal de La Calera"
$ws.Range("C73").Value = "Coquimbo"
$ws.Range("D73").Value = 44802
$ws.Range("E73").Value = 5
$ws.Range("F73").Value = 100112039
$ws.Range("G73").Value = "Ciboulette"
$ws.Range("H73").Value = "Sin especificar"
$ws.Range("I73").Value = "Primera"
$ws.Range("J73").Value = 120
$ws.Range("K73").Value = 1500
$ws.Range("L73").Value = 1500
$ws.Range("M73").Value = 1500
$ws.Range("N73").Value = "$/docena de atados"
$ws.Range("O73").Value = "Provincia de Quillota"
$ws.Range("P73").Value = 500
$ws.Range("Q73").Value = 3
$ws.Range("R73").Value = "Hortaliza"
